# Refined metadata to be additional tab
#
# 1. Add a new "metadata" worksheet after "data".
# 2. Update the F2:F14 "time_taken" timestamps on "data" to the re-run values.
# 3. Populate "metadata" with its header row + single data row, reusing the
#    same header/index cell formatting ("data"!B1 / "data"!A2) as the source
#    sheet via PasteSpecial(xlPasteFormats) so the shared style index is
#    reused instead of a new cellXfs entry being minted.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- update data!F2:F14 timestamps -----------------------------------------
$data.Range("F2").Value = "2021-10-05 14:18:59.406567"
$data.Range("F3").Value = "2021-10-05 14:18:59.406576"
$data.Range("F4").Value = "2021-10-05 14:18:59.406579"
$data.Range("F5").Value = "2021-10-05 14:18:59.406582"
$data.Range("F6").Value = "2021-10-05 14:18:59.406585"
$data.Range("F7").Value = "2021-10-05 14:18:59.406587"
$data.Range("F8").Value = "2021-10-05 14:18:59.406590"
$data.Range("F9").Value = "2021-10-05 14:18:59.406593"
$data.Range("F10").Value = "2021-10-05 14:18:59.406595"
$data.Range("F11").Value = "2021-10-05 14:18:59.406598"
$data.Range("F12").Value = "2021-10-05 14:18:59.406601"
$data.Range("F13").Value = "2021-10-05 14:18:59.406604"
$data.Range("F14").Value = "2021-10-05 14:18:59.406606"

# --- add the "metadata" sheet right after "data" ----------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse the bold/border/center style from data!B1:F1.
$data.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Index cell (A2) - reuse the style from data!A2.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

# Data row (B2:G2).
$meta.Range("B2").Value = "Additional findings health related - adults"
$meta.Range("C2").Value = 933

# "2.0" must stay text (not collapse to the number 2). Stage it through a
# scratch cell formatted as Text, then copy only the *value* across so the
# destination cell picks up the string without inheriting the scratch
# cell's NumberFormat override (keeps metadata!D2 on the default style).
$scratch = $meta.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2.0"
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$scratch.Clear()

$meta.Range("E2").Value = "2021-04-07T10:23:11.814083Z"
$meta.Range("F2").Value = "2021-10-05 14:18:59.402648"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/933/?format=json"

$data.Range("A1").Select()
